$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "43.279.52"
$ws.Range("E2").Value = "  -1.55%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.275.00"
$ws.Range("E3").Value = "  -1.81%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.12%  "

# Row 5 - Solana
$ws.Range("D5").Value = "'112.26"
$ws.Range("E5").Value = "  -2.16%  "

# Row 6 - BNB
$ws.Range("D6").Value = "'264.71"
$ws.Range("E6").Value = "  -2.14%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  -1.25%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.05%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "'0.608"
$ws.Range("E9").Value = "  -2.58%  "

# Row 10 - Avalanche
$ws.Range("D10").Value = "'47.82"
$ws.Range("E10").Value = "  +0.88%  "

# Row 11 - Dogecoin
$ws.Range("D11").Value = "'0.0931"
$ws.Range("E11").Value = "  -1.37%  "

# Row 12 - Polkadot
$ws.Range("D12").Value = "'8.81"
$ws.Range("E12").Value = "  -0.87%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +0.80%  "

# Row 14 - Chainlink
$ws.Range("D14").Value = "'15.48"
$ws.Range("E14").Value = "  -1.84%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "2.615.98"
$ws.Range("E15").Value = "  -1.79%  "

# Row 16 - Polygon
$ws.Range("D16").Value = "'0.855"
$ws.Range("E16").Value = "  -0.94%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.274.38"
$ws.Range("E17").Value = "  -1.96%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "43.172.15"
$ws.Range("E18").Value = "  -1.60%  "

# Row 19 - ShibaInu
$ws.Range("E19").Value = "  -2.49%  "

# Row 20 - Uniswap
$ws.Range("E20").Value = "  +1.63%  "

# Row 21 - Litecoin
$ws.Range("D21").Value = "'71.15"
$ws.Range("E21").Value = "  -2.26%  "

# Row 22 - ImmutableX
$ws.Range("D22").Value = "'2.51"
$ws.Range("E22").Value = "  +0.87%  "

# Row 23 - BitcoinCash
$ws.Range("D23").Value = "'231.68"
$ws.Range("E23").Value = "  -1.25%  "

# Row 24 - InternetComputer(DFINITY)
$ws.Range("D24").Value = "'9.69"
$ws.Range("E24").Value = "  +2.96%  "

# Row 25 - PancakeSwap
$ws.Range("D25").Value = "'2.88"
$ws.Range("E25").Value = "  -0.49%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  +0.28%  "

# Row 27 - Cosmos
$ws.Range("D27").Value = "'11.32"
$ws.Range("E27").Value = "  -1.47%  "

# Row 28 - LEO
$ws.Range("E28").Value = "  -1.06%  "

# Row 29 - InjectiveProtocol
$ws.Range("D29").Value = "'40.38"
$ws.Range("E29").Value = "  -5.39%  "

# Row 30 - Toncoin
$ws.Range("E30").Value = "  -1.70%  "

# Row 31 - WEMIXToken
$ws.Range("D31").Value = "'3.29"
$ws.Range("E31").Value = "  -4.15%  "

# Row 32 - Monero
$ws.Range("D32").Value = "'171.81"
$ws.Range("E32").Value = "  -3.47%  "

# Row 33 - EthereumClassic
$ws.Range("D33").Value = "'21.33"
$ws.Range("E33").Value = "  -2.95%  "

# Row 34 - Hedera
$ws.Range("D34").Value = "'0.0904"
$ws.Range("E34").Value = "  -3.37%  "

# Row 35 - Filecoin
$ws.Range("D35").Value = "'5.71"
$ws.Range("E35").Value = "  +2.09%  "

# Row 36 - Stellar
$ws.Range("E36").Value = "  +0.30%  "

# Row 37 - RenderToken
$ws.Range("D37").Value = "'4.65"
$ws.Range("E37").Value = "  -1.97%  "

# Row 38 - VeChain
$ws.Range("D38").Value = "'0.0352"
$ws.Range("E38").Value = "  -1.46%  "

# Row 39 - now NEARProtocol (was Kaspa)
$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D39").Value = "'3.83"
$ws.Range("E39").Value = "  -3.11%  "

# Row 40 - now Kaspa (was NEARProtocol)
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "'0.105"
$ws.Range("E40").Value = "  -6.29%  "

# Row 41 - LidoDAOToken
$ws.Range("D41").Value = "'2.60"
$ws.Range("E41").Value = "  +7.97%  "

# Row 42 - MultiversX
$ws.Range("D42").Value = "'76.53"
$ws.Range("E42").Value = "  +9.18%  "

# Row 43 - Celestia
$ws.Range("D43").Value = "'14.03"
$ws.Range("E43").Value = "  +10.88%  "

# Row 44 - Algorand
$ws.Range("D44").Value = "'0.236"
$ws.Range("E44").Value = "  -3.91%  "

# Row 45 - THORChain
$ws.Range("D45").Value = "'6.09"
$ws.Range("E45").Value = "  +2.67%  "

# Row 46 - FirstDigitalUSD
$ws.Range("E46").Value = "  -0.05%  "

# Row 47 - ARBITRUM
$ws.Range("D47").Value = "'1.37"
$ws.Range("E47").Value = "  -2.83%  "

# Row 48 - FraxShare
$ws.Range("D48").Value = "'8.61"
$ws.Range("E48").Value = "  -2.46%  "

# Row 49 - Cronos
$ws.Range("D49").Value = "'0.0993"
$ws.Range("E49").Value = "  -1.56%  "

# Row 50 - Aave
$ws.Range("D50").Value = "'101.56"
$ws.Range("E50").Value = "  +1.24%  "

# Row 51 - TrustWalletToken
$ws.Range("D51").Value = "'1.25"
$ws.Range("E51").Value = "  +1.70%  "
